$d = $word.ActiveDocument

# Locate the old email's local-part inside the "Contact Us" paragraph and
# replace it with the new one. Toggling a character property (Bold on/off)
# around the text assignment makes the engine keep this span as its own
# run instead of silently re-merging it with the identically formatted
# runs on either side - mirroring the three-run split seen in the diff.
$rng = $d.Content
$rng.Find.Execute("bartosz.broda98", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Text = "expirydateguard"
$rng.Bold = 1
$rng.Bold = 0
